# "Generate Report for Archive"
# The localization status of the "Ready for handoff" items moved on to
# "In Translation", so every cell holding that status text (on the
# Overview sheet as well as each per-locale sheet) is updated, and the
# Status column(s) are re-sized (narrower, since the new text is shorter
# than the old text) to match the refreshed report.

$wb = $excel.ActiveWorkbook

$oldStatus = "Ready for handoff"
$newStatus = "In Translation"

foreach ($ws in $wb.Worksheets) {
    $used = $ws.UsedRange
    foreach ($cell in $used.Cells) {
        # NOTE: keep the string literal on the left-hand side of -eq so
        # PowerShell does a string comparison instead of coercing the
        # literal into whatever type Value2 happens to be (e.g. cells
        # holding boolean TRUE/FALSE values).
        if ($oldStatus -eq $cell.Value2) {
            $cell.Value = $newStatus
        }
    }
}

# The Status columns auto-fit to the new (shorter) text: column width
# shrinks from ~17.22 characters down to ~13.41 characters.
$newColumnWidth = 12.5

$overview = $wb.Worksheets.Item("Overview")
$overview.Columns.Item(5).ColumnWidth = $newColumnWidth   # column E (zh-cn status)
$overview.Columns.Item(6).ColumnWidth = $newColumnWidth   # column F (de-de status)

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Columns.Item(3).ColumnWidth = $newColumnWidth        # column C (Status)

$dede = $wb.Worksheets.Item("de-de")
$dede.Columns.Item(3).ColumnWidth = $newColumnWidth        # column C (Status)
